$d = $word.ActiveDocument

# Locate the run that holds the full author list ("rank Gibson ... Ruttenberg" - the
# leading "F" of "Frank" lives in the previous run already).
$rng = $d.Content
$found = $rng.Find.Execute("rank Gibson, Allyson L. Lister, James Malone, Daniel Schober, Ryan R, Brinkman and Alan Ruttenberg", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Could not locate the author list run" }

$runStart = $rng.Start
$runEnd   = $rng.End

# Text as it should read after the edit: the run gets broken into three runs, and the
# stray comma after "Ryan R" becomes a period.
$part1    = "rank Gibson, Allyson L. Lister, James"
$part2New = " Malone, Daniel Schober, Ryan R."
$part2Old = " Malone, Daniel Schober, Ryan R,"
$part3    = " Brinkman and Alan Ruttenberg"

# Offsets (in the unedited text) of the two internal split points.
$split1    = $runStart + $part1.Length
$split2Old = $split1 + $part2Old.Length

# The run right after this one is just the closing "."; the run after that is a single
# space starting the next sentence. Remember both boundaries so they can be restored.
$periodEnd = $runEnd + 1
$regionEnd = $periodEnd + 1

# 1) Correct the text (comma -> period). Rewriting a range's text causes the host to
#    merge every adjacent, identically-formatted run in the paragraph into one run, so
#    fix the text first and rebuild the run boundaries afterwards.
$fullRange = $d.Range($runStart, $runEnd)
$fullRange.Text = $part1 + $part2New + $part3

# 2) Re-establish the run boundary before this text (between "...Courtot, F" and this
#    run), the two new internal split points, and the boundary after the closing period
#    (between it and the trailing space). Toggling a character property on and back off
#    forces the host to materialize a run break at that point without changing the
#    visible formatting. Each toggle range runs to a fixed point beyond every boundary
#    so none of them collapses to a zero-length range.
$boundaries = @($runStart, $split1, $split2Old, $runEnd, $periodEnd)
foreach ($pos in $boundaries) {
    $tail = $d.Range($pos, $regionEnd)
    $tail.Bold = 1
    $tail.Bold = 0
}
